# Updates the cryptocurrency price/volume table to reflect the latest
# scrape, matching the GitHub Actions "Updated cryptos list" commit.
# Also fixes the swapped order of EnergySwap / BabyDogeCoin rows (45-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = "D2"; Value = "29.599.69"; Text = 0 },
    @{ Ref = "E2"; Value = "  +2.35%  "; Text = 0 },
    @{ Ref = "D3"; Value = "1.859.00"; Text = 0 },
    @{ Ref = "E3"; Value = "  +1.46%  "; Text = 0 },
    @{ Ref = "D4"; Value = "0.9996"; Text = 1 },
    @{ Ref = "E4"; Value = "  +0.00%  "; Text = 0 },
    @{ Ref = "D5"; Value = "244.80"; Text = 1 },
    @{ Ref = "E5"; Value = "  +1.60%  "; Text = 0 },
    @{ Ref = "D6"; Value = "0.6945"; Text = 1 },
    @{ Ref = "E6"; Value = "  +1.11%  "; Text = 0 },
    @{ Ref = "D7"; Value = "1.000"; Text = 1 },
    @{ Ref = "E7"; Value = "  +0.04%  "; Text = 0 },
    @{ Ref = "D8"; Value = "0.07695"; Text = 1 },
    @{ Ref = "E8"; Value = "  +0.55%  "; Text = 0 },
    @{ Ref = "D9"; Value = "0.3059"; Text = 1 },
    @{ Ref = "E9"; Value = "  +0.38%  "; Text = 0 },
    @{ Ref = "D10"; Value = "23.71"; Text = 1 },
    @{ Ref = "E10"; Value = "  +0.61%  "; Text = 0 },
    @{ Ref = "D11"; Value = "0.07759"; Text = 1 },
    @{ Ref = "E11"; Value = "  -0.59%  "; Text = 0 },
    @{ Ref = "D12"; Value = "5.156"; Text = 1 },
    @{ Ref = "E12"; Value = "  +1.66%  "; Text = 0 },
    @{ Ref = "D13"; Value = "1.841.14"; Text = 0 },
    @{ Ref = "E13"; Value = "  +0.44%  "; Text = 0 },
    @{ Ref = "D14"; Value = "91.81"; Text = 1 },
    @{ Ref = "E14"; Value = "  +1.63%  "; Text = 0 },
    @{ Ref = "D15"; Value = "0.6913"; Text = 1 },
    @{ Ref = "E15"; Value = "  +2.14%  "; Text = 0 },
    @{ Ref = "D16"; Value = "6.569"; Text = 1 },
    @{ Ref = "E16"; Value = "  +1.74%  "; Text = 0 },
    @{ Ref = "D17"; Value = "29.599.91"; Text = 0 },
    @{ Ref = "E17"; Value = "  +2.35%  "; Text = 0 },
    @{ Ref = "D18"; Value = "0.000008292"; Text = 1 },
    @{ Ref = "E18"; Value = "  +0.23%  "; Text = 0 },
    @{ Ref = "D19"; Value = "2.103.42"; Text = 0 },
    @{ Ref = "E19"; Value = "  +1.22%  "; Text = 0 },
    @{ Ref = "D20"; Value = "240.33"; Text = 1 },
    @{ Ref = "E20"; Value = "  -0.92%  "; Text = 0 },
    @{ Ref = "E21"; Value = "  +0.78%  "; Text = 0 },
    @{ Ref = "D22"; Value = "1.000"; Text = 1 },
    @{ Ref = "E22"; Value = "  +0.05%  "; Text = 0 },
    @{ Ref = "D23"; Value = "7.596"; Text = 1 },
    @{ Ref = "E23"; Value = "  +2.30%  "; Text = 0 },
    @{ Ref = "E24"; Value = "  +0.05%  "; Text = 0 },
    @{ Ref = "D25"; Value = "0.1498"; Text = 1 },
    @{ Ref = "E25"; Value = "  +1.85%  "; Text = 0 },
    @{ Ref = "D26"; Value = "8.922"; Text = 1 },
    @{ Ref = "E26"; Value = "  +1.64%  "; Text = 0 },
    @{ Ref = "D27"; Value = "160.09"; Text = 1 },
    @{ Ref = "E27"; Value = "  -0.71%  "; Text = 0 },
    @{ Ref = "D28"; Value = "18.28"; Text = 1 },
    @{ Ref = "E28"; Value = "  +0.51%  "; Text = 0 },
    @{ Ref = "D29"; Value = "1.535"; Text = 1 },
    @{ Ref = "E29"; Value = "  +0.19%  "; Text = 0 },
    @{ Ref = "E30"; Value = "  +0.93%  "; Text = 0 },
    @{ Ref = "D31"; Value = "4.185"; Text = 1 },
    @{ Ref = "E31"; Value = "  +1.96%  "; Text = 0 },
    @{ Ref = "D32"; Value = "1.201"; Text = 1 },
    @{ Ref = "E32"; Value = "  +0.55%  "; Text = 0 },
    @{ Ref = "D33"; Value = "0.05093"; Text = 1 },
    @{ Ref = "E33"; Value = "  -0.44%  "; Text = 0 },
    @{ Ref = "D34"; Value = "0.7723"; Text = 1 },
    @{ Ref = "E34"; Value = "  +2.52%  "; Text = 0 },
    @{ Ref = "E35"; Value = "  +3.53%  "; Text = 0 },
    @{ Ref = "D36"; Value = "1.151"; Text = 1 },
    @{ Ref = "E36"; Value = "  +0.70%  "; Text = 0 },
    @{ Ref = "D37"; Value = "2.683"; Text = 1 },
    @{ Ref = "E37"; Value = "  +0.32%  "; Text = 0 },
    @{ Ref = "D38"; Value = "1.334.71"; Text = 0 },
    @{ Ref = "E38"; Value = "  +8.92%  "; Text = 0 },
    @{ Ref = "E39"; Value = "  +1.49%  "; Text = 0 },
    @{ Ref = "D40"; Value = "2.726"; Text = 1 },
    @{ Ref = "E40"; Value = "  +1.17%  "; Text = 0 },
    @{ Ref = "D41"; Value = "0.9646"; Text = 1 },
    @{ Ref = "E41"; Value = "  +4.71%  "; Text = 0 },
    @{ Ref = "D42"; Value = "106.75"; Text = 1 },
    @{ Ref = "E42"; Value = "  -1.41%  "; Text = 0 },
    @{ Ref = "D43"; Value = "5.791"; Text = 1 },
    @{ Ref = "E43"; Value = "  +3.04%  "; Text = 0 },
    @{ Ref = "D44"; Value = "0.9999"; Text = 1 },
    @{ Ref = "E44"; Value = "  +0.07%  "; Text = 0 },
    @{ Ref = "B45"; Value = "BabyDogeCoin"; Text = 0 },
    @{ Ref = "C45"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; Text = 0 },
    @{ Ref = "D45"; Value = "0.00000000126"; Text = 1 },
    @{ Ref = "E45"; Value = "  +3.66%  "; Text = 0 },
    @{ Ref = "B46"; Value = "EnergySwap"; Text = 0 },
    @{ Ref = "C46"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; Text = 0 },
    @{ Ref = "D46"; Value = "9.775"; Text = 1 },
    @{ Ref = "E46"; Value = "  +2.79%  "; Text = 0 },
    @{ Ref = "D47"; Value = "2.002.50"; Text = 0 },
    @{ Ref = "E48"; Value = "  +0.85%  "; Text = 0 },
    @{ Ref = "D49"; Value = "1.773"; Text = 1 },
    @{ Ref = "E49"; Value = "  +2.23%  "; Text = 0 },
    @{ Ref = "D50"; Value = "63.55"; Text = 1 },
    @{ Ref = "E50"; Value = "  -0.43%  "; Text = 0 },
    @{ Ref = "D51"; Value = "6.956"; Text = 1 },
    @{ Ref = "E51"; Value = "  +0.93%  "; Text = 0 }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.Text -eq 1) {
        # Force a text number format first so Excel keeps the exact
        # string (e.g. "1.000" / "0.9996") instead of re-parsing it as
        # a numeric value and stripping/rounding the trailing digits.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
